# Update country data and sharedStrings order per "Update countries & provincias Spain" commit.
# The underlying shared-string table was reshuffled (several country names
# swapped table positions) and the COVID-19 case figures were refreshed to
# the 18:22 snapshot (from 17:52). Net effect on the worksheet is a set of
# per-row updates: each affected row gets the country name now occupying
# that position plus that country's refreshed totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, País, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos críticos, Muertes hoy, Muertes
$updates = @(
    @(1, 'Datos actualizados a 6 de Abril de 2020 a las 18:22', $null, $null, $null, $null, $null, $null, $null),
    @(4, 'Estados Unidos', 349885, 13212, 19219, 320339, 8830, 711, 10327),
    @(6, 'Italia', 132547, 3599, 22837, 93187, 3898, 636, 16523),
    @(7, 'Alemania', 101089, 966, 28700, 70777, 3936, 28, 1612),
    @(17, 'Austria', 12267, 216, 3463, 8584, 250, 16, 220),
    @(19, 'Brasil', 11516, 262, 127, 10883, 296, 20, 506),
    @(29, 'India', 4693, 404, 346, 4218, 0, 11, 129),
    @(30, 'Dinamarca', 4681, 312, 1378, 3116, 144, 8, 118),
    @(34, 'Pakistan', 3662, 505, 259, 3351, 17, 5, 52),
    @(35, 'Filipinas', 3660, 414, 73, 3424, 1, 11, 163),
    @(48, 'Catar', 1832, 228, 131, 1697, 37, 0, 4),
    @(49, 'Republica Dominicana', 1828, 83, 33, 1709, 147, 4, 86),
    @(50, 'Grecia', 1755, 20, 269, 1407, 90, 6, 79),
    @(51, 'Sudafrica', 1655, 0, 95, 1549, 7, 0, 11),
    @(55, 'Argelia', 1423, 103, 90, 1160, 46, 21, 173),
    @(56, 'Singapur', 1375, 66, 344, 1025, 25, 0, 6),
    @(104, 'Montenegro', 233, 19, 1, 230, 4, 0, 2),
    @(105, 'Nigeria', 232, 0, 33, 194, 2, 0, 5),
    @(106, 'Senegal', 226, 4, 92, 132, 1, 0, 2),
    @(114, 'Consejo Danes para los Refugiados', 161, 7, 5, 138, 0, 0, 18),
    @(115, 'Venezuela', 159, 0, 52, 100, 6, 0, 7),
    @(116, 'Kenia', 158, 16, 4, 148, 2, 2, 6),
    @(176, 'Curazao', 13, 2, 5, 7, 0, 0, 1),
    @(178, 'Laos', 12, 1, 0, 12, 0, 0, 0),
    @(179, 'Sudan', 12, 0, 2, 8, 0, 0, 2),
    @(180, 'Seychelles', 11, 1, 0, 11, 0, 0, 0),
    @(181, 'Groenlandia', 11, 0, 3, 8, 0, 0, 0),
    @(184, 'Surinam', 10, 0, 0, 9, 0, 0, 1),
    @(185, 'Mozambique', 10, 0, 1, 9, 0, 0, 0)
)

foreach ($u in $updates) {
    $rowNum = $u[0]
    $ws.Cells.Item($rowNum, 1).Value = $u[1]
    for ($col = 2; $col -le 8; $col++) {
        $val = $u[$col]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $col).Value = $val
        }
    }
}
